# Update factsheets with text edits from COMM
#
# Converts numeric "count" cells to text-formatted cells (matching values),
# fixes the all-zero Nebraska counties in the County sheet to show
# percentage/currency placeholders instead of bare "0", and appends a
# "Total" row (row 84) to the County sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue($sheet, $row, $col, $val) {
    $c = $sheet.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Get-CellText($sheet, $row, $col) {
    $v = $sheet.Cells.Item($row, $col).Value2
    return "$v"
}

# ---------------------------------------------------------------------
# Sheet "Overall": A2 (count) becomes a text cell
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
$valA2 = Get-CellText $wsOverall 2 1
Set-TextValue $wsOverall 2 1 $valA2

# ---------------------------------------------------------------------
# Sheet "County": column B (counts) for rows 2-70 become text cells;
# the all-zero rows 71-83 get new placeholder text; a new Total row 84
# is appended.
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

for ($r = 2; $r -le 70; $r++) {
    $val = Get-CellText $wsCounty $r 2
    Set-TextValue $wsCounty $r 2 $val
}

for ($r = 71; $r -le 83; $r++) {
    Set-TextValue $wsCounty $r 2 "0.00%"
    Set-TextValue $wsCounty $r 3 "`$0"
    Set-TextValue $wsCounty $r 4 "0.00%"
    Set-TextValue $wsCounty $r 5 "0.00%"
    Set-TextValue $wsCounty $r 6 "0.00%"
}

Set-TextValue $wsCounty 84 1 "Total"
Set-TextValue $wsCounty 84 2 "769"
Set-TextValue $wsCounty 84 3 "`$986,876,099"
Set-TextValue $wsCounty 84 4 "7.78%"
Set-TextValue $wsCounty 84 5 "-11.28%"
Set-TextValue $wsCounty 84 6 "65.28%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": column B (counts) for rows 2-5
# ---------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")
for ($r = 2; $r -le 5; $r++) {
    $val = Get-CellText $wsCD $r 2
    Set-TextValue $wsCD $r 2 $val
}

# ---------------------------------------------------------------------
# Sheet "Size": column B (counts) for rows 2-8
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
for ($r = 2; $r -le 8; $r++) {
    $val = Get-CellText $wsSize $r 2
    Set-TextValue $wsSize $r 2 $val
}

# ---------------------------------------------------------------------
# Sheet "Subsector": column B (counts) for rows 2-13
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
for ($r = 2; $r -le 13; $r++) {
    $val = Get-CellText $wsSubsector $r 2
    Set-TextValue $wsSubsector $r 2 $val
}

Write-Host "Edits applied."
